$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7552.759
$ws.Range("J32").Value = 4861.5
$ws.Range("L32").Value = 4861.5
$ws.Range("N32").Value = -5513.5
$ws.Range("H43").Value = 6984.143
$ws.Range("I43").Value = 5631.3335
$ws.Range("J43").Value = 7998.75
$ws.Range("K43").Value = 5631.3335
$ws.Range("L43").Value = 7998.75
$ws.Range("M43").Value = -5562.3335
$ws.Range("N43").Value = -8136.75
$ws.Range("H53").Value = 131.16667
$ws.Range("I53").Value = 164.25
$ws.Range("K53").Value = 164.25
$ws.Range("M53").Value = 472.75
$ws.Range("H86").Value = 4612.923
$ws.Range("I86").Value = 4455.25
$ws.Range("J86").Value = 5138.5
$ws.Range("K86").Value = 4455.25
$ws.Range("L86").Value = 5138.5
$ws.Range("M86").Value = -3332.25
$ws.Range("N86").Value = -7384.5
$ws.Range("H89").Value = 4612.923
$ws.Range("I89").Value = 4455.25
$ws.Range("J89").Value = 5138.5
$ws.Range("K89").Value = 22276.25
$ws.Range("L89").Value = 25692.5
$ws.Range("M89").Value = -16660.25
$ws.Range("N89").Value = -36924.5
$ws.Range("H106").Value = 3356.6365
$ws.Range("I106").Value = 2992.3
$ws.Range("K106").Value = 2992.3
$ws.Range("M106").Value = -2361.3
$ws.Range("H130").Value = 120000
$ws.Range("J130").Value = 120000
$ws.Range("L130").Value = 120000
$ws.Range("H132").Value = 2135
$ws.Range("I132").Value = 2127.0952
$ws.Range("J132").Value = 2162.6667
$ws.Range("K132").Value = 6381.285600000001
$ws.Range("L132").Value = 6488.000100000001
$ws.Range("M132").Value = -3851.285600000001
$ws.Range("N132").Value = -11548.0001
$ws.Range("H138").Value = 3387.1353
$ws.Range("J138").Value = 3661.2666
$ws.Range("L138").Value = 10983.7998
$ws.Range("N138").Value = -21263.7998

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5967.6924
$ws.Range("I61").Value = 5598.273
$ws.Range("K61").Value = 5598.273
$ws.Range("M61").Value = -5386.273
$ws.Range("H88").Value = 2565.5715
$ws.Range("I88").Value = 1194
$ws.Range("J88").Value = 3594.25
$ws.Range("K88").Value = 1194
$ws.Range("L88").Value = 3594.25
$ws.Range("M88").Value = -788
$ws.Range("N88").Value = -4406.25
$ws.Range("H91").Value = 2565.5715
$ws.Range("I91").Value = 1194
$ws.Range("J91").Value = 3594.25
$ws.Range("K91").Value = 1194
$ws.Range("L91").Value = 3594.25
$ws.Range("M91").Value = 210
$ws.Range("N91").Value = -6402.25
$ws.Range("H136").Value = 5967.6924
$ws.Range("I136").Value = 5598.273
$ws.Range("K136").Value = 16794.819
$ws.Range("M136").Value = -14244.819

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1437.6666
$ws.Range("I20").Value = 1511.7858
$ws.Range("K20").Value = 1511.7858
$ws.Range("M20").Value = -1264.7858
$ws.Range("H134").Value = 2646.08
$ws.Range("I134").Value = 2652.1667
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 7956.500100000001
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -5421.500100000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 31623.2
$ws.Range("I16").Value = 22315.428
$ws.Range("K16").Value = 22315.428
$ws.Range("M16").Value = -22028.428
$ws.Range("H31").Value = 5717.7295
$ws.Range("I31").Value = 7056.5
$ws.Range("K31").Value = 7056.5
$ws.Range("M31").Value = -6761.5
$ws.Range("H34").Value = 5717.7295
$ws.Range("I34").Value = 7056.5
$ws.Range("K34").Value = 7056.5
$ws.Range("M34").Value = -6854.5
$ws.Range("H86").Value = 18526184
$ws.Range("I86").Value = 37043870
$ws.Range("J86").Value = 8497.888999999999
$ws.Range("K86").Value = 37043870
$ws.Range("L86").Value = 8497.888999999999
$ws.Range("M86").Value = -37042747
$ws.Range("N86").Value = -10743.889
$ws.Range("H89").Value = 18526184
$ws.Range("I89").Value = 37043870
$ws.Range("J89").Value = 8497.888999999999
$ws.Range("K89").Value = 185219350
$ws.Range("L89").Value = 42489.44499999999
$ws.Range("M89").Value = -185213734
$ws.Range("N89").Value = -53721.44499999999
$ws.Range("H107").Value = 2421.6
$ws.Range("I107").Value = 2649.5881
$ws.Range("K107").Value = 2649.5881
$ws.Range("M107").Value = -729.5880999999999
$ws.Range("H113").Value = 31623.2
$ws.Range("I113").Value = 22315.428
$ws.Range("K113").Value = 22315.428
$ws.Range("M113").Value = -20145.428
$ws.Range("H122").Value = 3002.9473
$ws.Range("I122").Value = 2989.8
$ws.Range("J122").Value = 3052.25
$ws.Range("K122").Value = 8969.400000000001
$ws.Range("L122").Value = 9156.75
$ws.Range("M122").Value = -6519.400000000001
$ws.Range("N122").Value = -14056.75
$ws.Range("H132").Value = 5744.7856
$ws.Range("I132").Value = 5772.5557
$ws.Range("K132").Value = 17317.6671
$ws.Range("M132").Value = -14787.6671
$ws.Range("H133").Value = 73764.664
$ws.Range("J133").Value = 77999
$ws.Range("L133").Value = 77999
$ws.Range("N133").Value = -83059
$ws.Range("H141").Value = 62352.082
$ws.Range("J141").Value = 64266.453
$ws.Range("L141").Value = 64266.453
$ws.Range("N141").Value = -74626.45300000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1088.3529
$ws.Range("I2").Value = 1242.3846
$ws.Range("J2").Value = 587.75
$ws.Range("K2").Value = 7454.3076
$ws.Range("L2").Value = 3526.5
$ws.Range("M2").Value = -7341.3076
$ws.Range("N2").Value = -3752.5
$ws.Range("H3").Value = 4459.5713
$ws.Range("I3").Value = 3203
$ws.Range("K3").Value = 9609
$ws.Range("M3").Value = -9497
$ws.Range("H12").Value = 367.82608
$ws.Range("J12").Value = 340.26666
$ws.Range("L12").Value = 1020.79998
$ws.Range("N12").Value = -1366.79998
$ws.Range("H16").Value = 5283
$ws.Range("I16").Value = 999.75
$ws.Range("J16").Value = 13849.5
$ws.Range("K16").Value = 2999.25
$ws.Range("L16").Value = 41548.5
$ws.Range("M16").Value = -2826.25
$ws.Range("N16").Value = -41894.5
$ws.Range("H131").Value = 2296.3333
$ws.Range("I131").Value = 1148.5454
$ws.Range("J131").Value = 2703.6128
$ws.Range("K131").Value = 3445.6362
$ws.Range("L131").Value = 8110.8384
$ws.Range("M131").Value = 1594.3638
$ws.Range("N131").Value = -18190.8384

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 15832.333
$ws.Range("I5").Value = 13999
$ws.Range("J5").Value = 24999
$ws.Range("K5").Value = 13999
$ws.Range("L5").Value = 24999
$ws.Range("M5").Value = -13887
$ws.Range("H70").Value = 5623
$ws.Range("I70").Value = 5330.6665
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 5330.6665
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -5060.6665
$ws.Range("N70").Value = -7040
$ws.Range("H73").Value = 5623
$ws.Range("I73").Value = 5330.6665
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 5330.6665
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = -4394.6665
$ws.Range("N73").Value = -8372
$ws.Range("H122").Value = 4708.4
$ws.Range("I122").Value = 4033.5386
$ws.Range("K122").Value = 12100.6158
$ws.Range("M122").Value = -9650.6158
$ws.Range("H132").Value = 4193.8857
$ws.Range("I132").Value = 4179.1377
$ws.Range("J132").Value = 4265.1665
$ws.Range("K132").Value = 12537.4131
$ws.Range("L132").Value = 12795.4995
$ws.Range("M132").Value = -10007.4131
$ws.Range("N132").Value = -17855.4995

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 492.1875
$ws.Range("I16").Value = 464.86667
$ws.Range("K16").Value = 464.86667
$ws.Range("M16").Value = -294.86667
$ws.Range("H46").Value = 5363.091
$ws.Range("J46").Value = 5124.5
$ws.Range("L46").Value = 5124.5
$ws.Range("N46").Value = -5500.5
$ws.Range("H141").Value = 82480.09
$ws.Range("J141").Value = 82480.09
$ws.Range("L141").Value = 82480.09
$ws.Range("N141").Value = -92840.09

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3065.5
$ws.Range("I132").Value = 2688.8064
$ws.Range("K132").Value = 8066.4192
$ws.Range("M132").Value = -5536.4192
$ws.Range("H135").Value = 72207.53999999999
$ws.Range("J135").Value = 72207.53999999999
$ws.Range("L135").Value = 72207.53999999999
$ws.Range("N135").Value = -82347.53999999999
$ws.Range("H136").Value = 1910.9697
$ws.Range("I136").Value = 1357.6552
$ws.Range("J136").Value = 5922.5
$ws.Range("K136").Value = 4072.9656
$ws.Range("L136").Value = 17767.5
$ws.Range("M136").Value = -1522.9656
$ws.Range("N136").Value = -22867.5

# --- New cells added (previously empty) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N130").Value = -130040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N134").Value = -12570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N5").Value = -25223
